$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("jan")
$ws2 = $wb.Worksheets.Item("feb")

# Add column G values on "jan" sheet (rows 15-20), mirroring column H
$ws1.Range("G15").Value = 4
$ws1.Range("G16").Value = 3
$ws1.Range("G17").Value = 4
$ws1.Range("G18").Value = 4
$ws1.Range("G19").Value = 4
$ws1.Range("G20").Value = 3

# Add column G values on "feb" sheet (rows 7-12), mirroring column H
$ws2.Range("G7").Value = 3
$ws2.Range("G8").Value = 5
$ws2.Range("G9").Value = 4
$ws2.Range("G10").Value = 4
$ws2.Range("G11").Value = 4
$ws2.Range("G12").Value = 4

# Update active sheet / selections: "feb" loses tab selection (was active),
# "jan" becomes the active tab with selection G20; "feb" selection becomes I14
$ws2.Activate()
$ws2.Range("I14").Select()

$ws1.Activate()
$ws1.Range("G20").Select()
